$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# The whole diagram lives inside one big top-level group ("Group 73").
$grp = $s.Shapes.Item(1)

# 1) Give the two "drop-down" arrow connectors (into the P-MPNN Input
#    Embedder box and into the Graph Reduction box) an explicit grey
#    outline colour.
$arrow1 = $grp.GroupItems.Item("Straight Arrow Connector 69")
$arrow1.Line.ForeColor.RGB = 8355711

$arrow2 = $grp.GroupItems.Item("Straight Arrow Connector 71")
$arrow2.Line.ForeColor.RGB = 8355711

# 2) Rename/resize the "Symmetry Remodeling" caption to "Graph Reduction"
#    (narrower text needs a narrower auto-fit text box).
$caption = $grp.GroupItems.Item("TextBox 72")
$caption.TextFrame.TextRange.Text = "Graph Reduction"
$caption.Width = 104.40976377952757

# 3) Add the new (mostly off-slide / incidental) grey-outlined textbox
#    that sits directly on the slide, alongside the group.
$tb = $s.Shapes.AddTextbox(1, 10, 20, 30, 40)
$tb.TextFrame.WordWrap = 0
$tb.TextFrame.AutoSize = 1
$tb.Fill.Visible = 0
$tb.Line.ForeColor.RGB = 8355711
$tb.Left = 27.391339302062992
$tb.Top = -87.65228652954102
$tb.Width = 14.545748233795168
